# Crackteck-modules.xlsx -- "added auth and leads apis"
#
# The sheet is a manual QA/status tracker. Column E holds a colored status
# tag ("Done" / "Pending" / "Not Done" / "Some Changes are pending") per
# feature row, rendered through a handful of reusable cell styles (fill +
# border + alignment). We flip several rows' status by copying the
# formatting (fill/border/alignment) from a donor cell that already carries
# the target style, then writing the matching status text so the shared
# string lines up with the style. We also bump one row's height, retarget
# the description of the "Low Stock Reports" row, and move the sheet's
# active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Row 31: grew taller (e.g. to fit wrapped/longer text after the edit) ---
$ws.Rows.Item(31).RowHeight = 20.25

# --- E35: "Pending" -> "Done" (auth API finished) ---
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E35").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E35").Value = $ws.Range("E7").Value()

# --- E36: "Some Changes are pending" -> "Done" ---
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E36").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E36").Value = $ws.Range("E7").Value()

# --- E59: "Not Done" -> "Some Changes are pending" (leads API in progress) ---
$ws.Range("E30").Copy() | Out-Null
$ws.Range("E59").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E59").Value = $ws.Range("E30").Value()

# --- E65: "Some Changes are pending" -> "Done" ---
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E65").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E65").Value = $ws.Range("E7").Value()

# --- E66: "Not Done" -> "Done" ---
$ws.Range("E7").Copy() | Out-Null
$ws.Range("E66").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("E66").Value = $ws.Range("E7").Value()

$excel.CutCopyMode = $false

# --- C80: describe the Low Stock Reports row using its own title text ---
$ws.Range("C80").Value = $ws.Range("B80").Value()

# --- move the active selection the way the author left it ---
$ws.Range("E59").Select() | Out-Null
